$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = 9
$ws.Range("B9").Value = "cunt"

$ws.Range("A10").Value = 10
$ws.Range("B10").Value = "hell"

$ws.Range("A11").Value = 11
$ws.Range("B11").Value = "bloody hell"
